$wb = $excel.ActiveWorkbook

# Add the new "LM2735" worksheet after the last existing sheet (ADC)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "LM2735"

# Column width
$ws.Columns.Item(2).ColumnWidth = 5.75

# Seed shared strings in the exact order the original author created them:
# R1, Vout, R2, then the two section titles.
$ws.Range("C5").Value = "R1"
$ws.Range("C4").Value = "Vout"
$ws.Range("C6").Value = "R2"
$ws.Range("B3").Value = "Vout에 대한 R2 계산"
$ws.Range("B8").Value = "R1,R2에 대한 Vout 계산"

# Section 1: compute R2 given Vout and R1
$ws.Range("B3").Font.Bold = $true

$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 5

$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 10

$ws.Range("D6").Formula = "=(D4/1.255-1)*D5"
$ws.Range("E6").Formula = "=(E4/1.255-1)*E5"
$ws.Range("C4:E6").Borders.LineStyle = 1

# Section 2: compute Vout given R1 and R2
$ws.Range("B8").Font.Bold = $true

$ws.Range("C9").Value = "R1"
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 10

$ws.Range("C10").Value = "R2"
$ws.Range("D10").Value = 56
$ws.Range("E10").Value = 30

$ws.Range("C11").Value = "Vout"
$ws.Range("D11").Formula = "=(D10/D9+1)*1.255"
$ws.Range("E11").Formula = "=(E10/E9+1)*1.255"
$ws.Range("C9:E11").Borders.LineStyle = 1

$ws.Range("D6:E6").NumberFormat = "0.0"
$ws.Range("D11:E11").NumberFormat = "0.0"

$ws.Range("H10").Select()

$ws.Activate()
